$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.604.75'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').Value = '1.666.42'
$ws.Range('E3').Value = '  -3.47%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.512'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.16%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.73'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0623'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.57%  '
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('D12').Value = '1.901.54'
$ws.Range('E12').Value = '  -3.52%  '
$ws.Range('D13').Value = '1.655.45'
$ws.Range('E13').Value = '  -4.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.15'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.564'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '27.619.63'
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '242.73'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('E19').Value = '  -3.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.43%  '
$ws.Range('E24').Value = '  -4.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.21'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.12%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  -2.24%  '
$ws.Range('E30').Value = '  +3.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0502'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.58%  '
$ws.Range('D33').Value = '1.473.01'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('E34').Value = '  -4.67%  '
$ws.Range('E35').Value = '  -5.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.936'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.38'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0172'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.574'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '69.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('E41').Value = '  -4.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.08%  '
$ws.Range('E44').Value = '  -3.74%  '
$ws.Range('D45').Value = '1.809.72'
$ws.Range('E45').Value = '  -3.43%  '
$ws.Range('E46').Value = '  -1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.71'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '89.44'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('E49').Value = '  -3.89%  '
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.88'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.88%  '
